# Auto-generated Excel COM-interop script
# Applies updated Leve profit calculations across multiple sheets
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 1681.5
$ws.Range("I80").Value = 1899.8
$ws.Range("J80").Value = 1317.6666
$ws.Range("K80").Value = 5699.4
$ws.Range("L80").Value = 3952.9998
$ws.Range("M80").Value = -4701.4
$ws.Range("N80").Value = -5948.9998
# Row 83
$ws.Range("H83").Value = 1681.5
$ws.Range("I83").Value = 1899.8
$ws.Range("J83").Value = 1317.6666
$ws.Range("K83").Value = 17098.2
$ws.Range("L83").Value = 11858.9994
$ws.Range("M83").Value = -12106.2
$ws.Range("N83").Value = -21842.9994
# Row 116
$ws.Range("H116").Value = 5594.3076
$ws.Range("I116").Value = 7537.8237
$ws.Range("J116").Value = 1923.2222
$ws.Range("K116").Value = 7537.8237
$ws.Range("L116").Value = 1923.2222
$ws.Range("M116").Value = -4095.8237
$ws.Range("N116").Value = -8807.2222
# Row 132
$ws.Range("H132").Value = 861.9804
$ws.Range("I132").Value = 807.5333000000001
$ws.Range("J132").Value = 1270.3334
$ws.Range("K132").Value = 2422.5999
$ws.Range("L132").Value = 3811.0002
$ws.Range("M132").Value = 107.4000999999998
$ws.Range("N132").Value = -8871.0002

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4901.95
$ws.Range("I32").Value = 3793.4065
$ws.Range("J32").Value = 16110.556
$ws.Range("K32").Value = 3793.4065
$ws.Range("L32").Value = 16110.556
$ws.Range("M32").Value = -3506.4065
$ws.Range("N32").Value = -16684.556
# Row 61
$ws.Range("H61").Value = 21281.295
$ws.Range("I61").Value = 23383.445
$ws.Range("J61").Value = 11821.625
$ws.Range("K61").Value = 23383.445
$ws.Range("L61").Value = 11821.625
$ws.Range("M61").Value = -23171.445
$ws.Range("N61").Value = -12245.625
# Row 110
$ws.Range("H110").Value = 1000
$ws.Range("I110").Value = 1000
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1000
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1045
# Row 122
$ws.Range("H122").Value = 1244.7222
$ws.Range("I122").Value = 1170.8823
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 3512.6469
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -1062.6469
$ws.Range("N122").Value = -12400
# Row 136
$ws.Range("H136").Value = 21281.295
$ws.Range("I136").Value = 23383.445
$ws.Range("J136").Value = 11821.625
$ws.Range("K136").Value = 70150.33499999999
$ws.Range("L136").Value = 35464.875
$ws.Range("M136").Value = -67600.33499999999
$ws.Range("N136").Value = -40564.875

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 339400
$ws.Range("I86").Value = 7280
$ws.Range("J86").Value = 2000000
$ws.Range("K86").Value = 7280
$ws.Range("L86").Value = 2000000
$ws.Range("M86").Value = -6157
$ws.Range("N86").Value = -2002246
# Row 89
$ws.Range("H89").Value = 339400
$ws.Range("I89").Value = 7280
$ws.Range("J89").Value = 2000000
$ws.Range("K89").Value = 36400
$ws.Range("L89").Value = 10000000
$ws.Range("M89").Value = -30784
$ws.Range("N89").Value = -10011232
# Row 99
$ws.Range("H99").Value = 1335.1818
$ws.Range("I99").Value = 1185.8334
$ws.Range("J99").Value = 1514.4
$ws.Range("K99").Value = 1185.8334
$ws.Range("L99").Value = 1514.4
$ws.Range("M99").Value = 312.1666
$ws.Range("N99").Value = -4510.4
# Row 107
$ws.Range("H107").Value = 900
$ws.Range("I107").Value = 820
$ws.Range("J107").Value = 1100
$ws.Range("K107").Value = 820
$ws.Range("L107").Value = 1100
$ws.Range("M107").Value = 1100
$ws.Range("N107").Value = -4940
# Row 134
$ws.Range("H134").Value = 4024.7441
$ws.Range("I134").Value = 4136.436
$ws.Range("J134").Value = 2935.75
$ws.Range("K134").Value = 12409.308
$ws.Range("L134").Value = 8807.25
$ws.Range("M134").Value = -9874.307999999999
$ws.Range("N134").Value = -13877.25

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Range("H134").Value = 3220.12
$ws.Range("I134").Value = 2934.9565
$ws.Range("J134").Value = 6499.5
$ws.Range("K134").Value = 8804.869499999999
$ws.Range("L134").Value = 19498.5
$ws.Range("M134").Value = -6269.869499999999
$ws.Range("N134").Value = -24568.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 1960
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 2075
$ws.Range("K80").Value = 4500
$ws.Range("L80").Value = 6225
$ws.Range("M80").Value = -3564
$ws.Range("N80").Value = -8097
# Row 83
$ws.Range("H83").Value = 1960
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 2075
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 18675
$ws.Range("M83").Value = -8820
$ws.Range("N83").Value = -28035
# Row 131
$ws.Range("H131").Value = 18530.8
$ws.Range("I131").Value = 1030
$ws.Range("J131").Value = 18928.545
$ws.Range("K131").Value = 3090
$ws.Range("L131").Value = 56785.63499999999
$ws.Range("M131").Value = 1950
$ws.Range("N131").Value = -66865.63499999999
# Row 137
$ws.Range("H137").Value = 4788.6313
$ws.Range("I137").Value = 2121.889
$ws.Range("J137").Value = 7188.7
$ws.Range("K137").Value = 6365.667
$ws.Range("L137").Value = 21566.1
$ws.Range("M137").Value = -1265.667
$ws.Range("N137").Value = -31766.1

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2431.7646
$ws.Range("I102").Value = 2481.4285
$ws.Range("J102").Value = 2200
$ws.Range("K102").Value = 2481.4285
$ws.Range("L102").Value = 2200
$ws.Range("M102").Value = -859.4285
$ws.Range("N102").Value = -5444

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1723.3
$ws.Range("I46").Value = 1273.5
$ws.Range("J46").Value = 2023.1666
$ws.Range("K46").Value = 1273.5
$ws.Range("L46").Value = 2023.1666
$ws.Range("M46").Value = -1085.5
$ws.Range("N46").Value = -2399.1666
# Row 55
$ws.Range("H55").Value = 460.6
$ws.Range("I55").Value = 399.5
$ws.Range("J55").Value = 521.7
$ws.Range("K55").Value = 399.5
$ws.Range("L55").Value = 521.7
$ws.Range("M55").Value = -226.5
$ws.Range("N55").Value = -867.7
# Row 61
$ws.Range("H61").Value = 4959.8
$ws.Range("I61").Value = 4999
$ws.Range("J61").Value = 4950
$ws.Range("K61").Value = 4999
$ws.Range("L61").Value = 4950
$ws.Range("M61").Value = -4797
$ws.Range("N61").Value = -5354
# Row 100
$ws.Range("H100").Value = 1610.5
$ws.Range("I100").Value = 1583.4286
$ws.Range("J100").Value = 1800
$ws.Range("K100").Value = 1583.4286
$ws.Range("L100").Value = 1800
$ws.Range("M100").Value = -1042.4286
$ws.Range("N100").Value = -2882
# Row 113
$ws.Range("H113").Value = 4959.8
$ws.Range("I113").Value = 4999
$ws.Range("J113").Value = 4950
$ws.Range("K113").Value = 4999
$ws.Range("L113").Value = 4950
$ws.Range("M113").Value = -2829
$ws.Range("N113").Value = -9290
# Row 132
$ws.Range("H132").Value = 2103.1904
$ws.Range("I132").Value = 1648.8334
$ws.Range("J132").Value = 2443.9583
$ws.Range("K132").Value = 4946.5002
$ws.Range("L132").Value = 7331.874899999999
$ws.Range("M132").Value = -2416.5002
$ws.Range("N132").Value = -12391.8749

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1197
$ws.Range("I100").Value = 1036.4
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 2072.8
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -1531.8
$ws.Range("N100").Value = -5082
# Row 104
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").ClearContents()
$ws.Range("N104").Value = 0
# Row 105
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value = 0
# Row 106
$ws.Range("H106").Value = 45000
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 45000
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 45000
$ws.Range("N106").Value = -47524
# Row 107
$ws.Range("H107").Value = 712.6842
$ws.Range("I107").Value = 571.3125
$ws.Range("J107").Value = 1466.6666
$ws.Range("K107").Value = 1713.9375
$ws.Range("L107").Value = 4399.9998
$ws.Range("M107").Value = 206.0625
$ws.Range("N107").Value = -8239.9998
# Row 136
$ws.Range("H136").Value = 1830.5186
$ws.Range("I136").Value = 1092.7858
$ws.Range("J136").Value = 2625
$ws.Range("K136").Value = 3278.3574
$ws.Range("L136").Value = 7875
$ws.Range("M136").Value = -728.3574000000003
$ws.Range("N136").Value = -12975
